$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right before "总计".
#    We duplicate "总计" first (so the new sheet inherits the exact same
#    header/row styling, margins, etc.), then rename it and overwrite its
#    contents with the 2022-Q1 fund-holding data.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Copy($total)
$q1 = $wb.Worksheets.Item("总计 (2)")
$q1.Name = "2022-Q1"

# Make sure the header row (currently only B:D) also covers E:H with the
# same style as the existing header cells.
$q1.Range("B1:D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("B2").Value = "'006923"
$q1.Range("C2").Value = "前海开源沪港深非周期性行业股票A"
$q1.Range("D2").Value = "'0.54"
$q1.Range("E2").Value = "'93.77"
$q1.Range("F2").Value = "'5.38"
$q1.Range("G2").Value = "'0.0291"
$q1.Range("H2").Value = 6

$q1.Range("B3").Value = "'006924"
$q1.Range("C3").Value = "前海开源沪港深非周期性行业股票C"
$q1.Range("D3").Value = "'0.22"
$q1.Range("E3").Value = "'93.77"
$q1.Range("F3").Value = "'5.38"
$q1.Range("G3").Value = "'0.0118"
$q1.Range("H3").Value = 6

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: add a new top data row for 2022-Q1 and
#    shift the existing 2021-Q4 / 2021-Q3 rows down by one.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")

$b2 = $ws.Range("B2").Value()
$c2 = $ws.Range("C2").Value()
$d2 = $ws.Range("D2").Value()
$b3 = $ws.Range("B3").Value()
$c3 = $ws.Range("C3").Value()
$d3 = $ws.Range("D3").Value()

# Row 4 <- old row 3 (2021-Q3)
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = $b3
$ws.Range("C4").Value = $c3
$ws.Range("D4").Value = $d3

# Row 3 <- old row 2 (2021-Q4)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = $b2
$ws.Range("C3").Value = $c2
$ws.Range("D3").Value = $d2

# Row 2 <- new 2022-Q1 summary
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0.04

# ---------------------------------------------------------------------
# 3) Restore the originally active sheet/tab (copying "总计" above makes
#    the new sheet active as a side effect).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
